## Grade update as of 10:07 AM Friday, January 3, 2025 (GMT+8)
## Time in 40 (Pob.)
##
## Updates individual student grade cells on sheets "B1" and "B3" and
## moves the active selection on "B3" to F50 (the cell last selected when
## the grades were saved).

$wb = $excel.ActiveWorkbook

# --- Sheet "B1" (1st sheet) ------------------------------------------------
$ws1 = $wb.Worksheets.Item("B1")

# FLORES,LINDELLE JOYCE ROCETE - NewtonRhapson score correction
$ws1.Range("G8").Value = 95

# MARGIN,NICOLE MACARAT - FINALS recomputed
$ws1.Range("I18").Value = 66.67

# MEDIANO,ALDREI JOSHUA SALE - Bisection & Bairstows scores recorded (were missing/0)
$ws1.Range("F19").Value = 100
$ws1.Range("H19").Value = 60

# MONTEBON,CELINE FE WASLO - FINALS recomputed
$ws1.Range("I20").Value = 66.67

# OLIFERNES,CHARLES ANTHONY BUNOTAN - FINALS recomputed
$ws1.Range("I21").Value = 66.67

# WONG,ALEX FERDIE MUÑASQUE - Arithmetic Series score recorded & FINALS recomputed
$ws1.Range("C28").Value = 100
$ws1.Range("I28").Value = 66.67

# --- Sheet "B3" (3rd sheet) -------------------------------------------------
$ws3 = $wb.Worksheets.Item("B3")

# MITMUG,KHADIJAH B - Arithmetic Series & Euler scores recorded (were missing/0)
$ws3.Range("C23").Value = 100
$ws3.Range("D23").Value = 70

# Leave the active sheet / selection where it was when last saved
$ws3.Activate() | Out-Null
$ws3.Range("F50").Select() | Out-Null
